$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for rows 2-340.
# This update bumps every one of those serials by exactly one day
# (45205 -> 45206, i.e. 2023-10-06 -> 2023-10-07).
for ($r = 2; $r -le 340; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    if ($v -eq 45205) {
        $cell.Value = 45206
    }
}
